$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D/E price & volume updates ---
$ws.Range('D2').Value = '26.779.39'
$ws.Range('E2').Value = '  +7.54%  '
$ws.Range('D3').Value = '1.743.07'
$ws.Range('E3').Value = '  +4.22%  '
$ws.Range('D4').Value = '''1.005'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '''334.06'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').Value = '''0.9950'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = '''0.3730'
$ws.Range('E7').Value = '  +2.43%  '
$ws.Range('D8').Value = '''48.65'
$ws.Range('E8').Value = '  +3.82%  '
$ws.Range('D9').Value = '''0.3387'
$ws.Range('E9').Value = '  +4.70%  '
$ws.Range('D10').Value = '''1.186'
$ws.Range('E10').Value = '  +4.02%  '
$ws.Range('D11').Value = '''0.07507'
$ws.Range('E11').Value = '  +6.39%  '
$ws.Range('D12').Value = '''0.9922'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').Value = '''6.373'
$ws.Range('E13').Value = '  +4.87%  '
$ws.Range('D14').Value = '''20.43'
$ws.Range('E14').Value = '  +4.10%  '
$ws.Range('D15').Value = '''7.046'
$ws.Range('E15').Value = '  +6.42%  '
$ws.Range('D16').Value = '1.750.12'
$ws.Range('E16').Value = '  +4.87%  '
$ws.Range('D17').Value = '''0.00001083'
$ws.Range('E17').Value = '  +3.62%  '
$ws.Range('D18').Value = '''0.06696'
$ws.Range('E18').Value = '  +2.29%  '
$ws.Range('D19').Value = '''82.86'
$ws.Range('E19').Value = '  +5.26%  '
$ws.Range('D20').Value = '''0.9969'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '''16.71'
$ws.Range('E21').Value = '  +5.42%  '
$ws.Range('D22').Value = '''6.236'
$ws.Range('E22').Value = '  +5.50%  '
$ws.Range('D23').Value = '''12.86'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '26.881.20'
$ws.Range('E24').Value = '  +7.78%  '
$ws.Range('D25').Value = '''2.445'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '''1.476'
$ws.Range('E26').Value = '  +25.63%  '
$ws.Range('D27').Value = '''2.422'
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D28').Value = '''152.01'
$ws.Range('E28').Value = '  +2.54%  '
$ws.Range('D29').Value = '''19.66'
$ws.Range('E29').Value = '  +4.98%  '
$ws.Range('D30').Value = '1.944.44'
$ws.Range('E30').Value = '  +5.14%  '
$ws.Range('D31').Value = '''133.12'
$ws.Range('E31').Value = '  +5.91%  '
$ws.Range('D32').Value = '''4.106'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = '''6.027'
$ws.Range('E33').Value = '  +4.00%  '
$ws.Range('D34').Value = '''0.08582'
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('D35').Value = '''1.687'
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('D36').Value = '''12.92'
$ws.Range('E36').Value = '  +5.35%  '
$ws.Range('D37').Value = '''5.435'
$ws.Range('E37').Value = '  +5.46%  '
$ws.Range('D38').Value = '''0.02338'
$ws.Range('E38').Value = '  +4.65%  '
$ws.Range('D39').Value = '''0.06272'
$ws.Range('E39').Value = '  +3.94%  '
$ws.Range('D40').Value = '''0.2169'
$ws.Range('E40').Value = '  +4.07%  '
$ws.Range('D41').Value = '''8.503'
$ws.Range('E41').Value = '  +3.36%  '
$ws.Range('D44').Value = '''14.34'
$ws.Range('E44').Value = '  +4.31%  '
$ws.Range('D45').Value = '''0.9949'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = '''0.6224'
$ws.Range('E46').Value = '  +8.69%  '
$ws.Range('D47').Value = '''3.920'
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('E50').Value = '  +3.12%  '
$ws.Range('D51').Value = '''78.27'
$ws.Range('E51').Value = '  +5.24%  '
# --- Row 42/43: coin identity swap (TheSandbox <-> TrustWalletToken) ---
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'1.218"
$ws.Range('E42').Value = '  -0.89%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = "'0.6299"
$ws.Range('E43').Value = '  +5.92%  '

# --- Row 48/49: coin identity swap (NEARProtocol <-> Quant) ---
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'129.17"
$ws.Range('E48').Value = '  +3.63%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'2.077"
$ws.Range('E49').Value = '  +5.94%  '
